$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "tables;/tables/list-subcategories;/tables/add-subcategories"
$ws.Range("C6").Value = "Anagrafica Sotto Gruppi"

$ws.Range("B6").Select()
